# Apply updated crypto price/volume values to the worksheet.
# D-column "price" cells can look numeric (e.g. "605.16", "0.999") but must
# stay as text, matching the original inlineStr cells, so we temporarily force
# a text number format before assigning the value, then restore the style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.135.78'
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = '  -0.17%  '
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.554.31'
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = '  +1.76%  '
$ws.Range("E4").Value = '  +0.00%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.16'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +0.06%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.15'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -0.07%  '
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.553.39'
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("E11").Value = '  -3.35%  '
$ws.Range("E12").Value = '  -0.26%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.155.48'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +1.75%  '
$ws.Range("E14").Value = '  +0.86%  '
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.99'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -1.27%  '
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.547.58'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  +1.69%  '
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.160.35'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("E19").Value = '  +5.44%  '
$ws.Range("E20").Value = '  +0.43%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.67'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -1.30%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '429.87'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").Value = '  +2.37%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.91'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  +2.50%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.696.31'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  +1.90%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("E28").Value = '  +0.66%  '
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.08'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -2.28%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.82'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -1.22%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = $style
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.550.81'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  +2.05%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.41'
$ws.Range("D33").Style = $style
$ws.Range("E34").Value = '  -2.22%  '
$ws.Range("E35").Value = '  -8.89%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.78'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  -1.38%  '
$ws.Range("E39").Value = '  -1.53%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '173.69'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  +2.16%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0845'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -1.90%  '
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.18'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  +0.15%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.890'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("E44").Value = '  +0.84%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.03'
$ws.Range("D45").Style = $style
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("E47").Value = '  -1.05%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.94'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("E51").Value = '  +2.84%  '
